$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Analog")

# A new point-name string ("HWE-IINYOKERN C TGT") was inserted into the
# shared-string table right after "SEL-2020 COM" (row 33's label) but the
# existing rows below it were not re-pointed to account for the shift, so
# each of them now shows the label that used to belong to the row above:
$ws.Range("B34").Value = "HWE-IINYOKERN C TGT"
$ws.Range("B35").Value = "lcoso HWE A TGT"
$ws.Range("B36").Value = "lcoso HWE B TGT"
$ws.Range("B37").Value = "lcoso HWE GRD TGT"

# New row 38: point 36, whose label ended up being the old row 37's text
# ("SILVERPEAK C C TCT") because of the same shift.
$ws.Range("A38").Value = "'36"
$ws.Range("A38").ClearFormats()
$ws.Range("B38").Value = "SILVERPEAK C C TCT"
